$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "51.224.20"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.86%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.958.01"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.57%  "

# Row 4
$ws.Range("E4").Value = "  +0.06%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "382.02"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.43%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "102.66"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.43%  "

# Row 7
$ws.Range("E7").Value = "  -0.36%  "

# Row 8
$ws.Range("E8").Value = "  +0.04%  "

# Row 9
$ws.Range("E9").Value = "  -1.08%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "36.62"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.71%  "

# Row 11
$ws.Range("E11").Value = "  -0.26%  "

# Row 12
$ws.Range("E12").Value = "  -0.28%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "3.423.83"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.76%  "

# Row 14
$ws.Range("E14").Value = "  -2.91%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.43"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.21%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.949.50"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.43%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.990"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +3.81%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "51.149.25"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.86%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.21"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -5.91%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.14"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -3.52%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.57"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -4.63%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.0₃0955"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.09%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "68.47"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.12%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "262.53"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.24%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.91"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +3.14%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "8.42"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +13.35%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.74"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +4.35%  "

# Row 28
$ws.Range("E28").Value = "  +0.88%  "

# Row 29
$ws.Range("E29").Value = "  -0.04%  "

# Row 30
$ws.Range("E30").Value = "  +8.72%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "25.72"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.83%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "9.81"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.80%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0458"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +5.29%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "33.85"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.65%  "

# Row 35
$ws.Range("B35").Value = "Toncoin"
$ws.Range("C35").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.05"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.75%  "

# Row 36
$ws.Range("B36").Value = "OKB"
$ws.Range("C36").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "50.37"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -3.64%  "

# Row 37
$ws.Range("E37").Value = "  +0.04%  "

# Row 38
$ws.Range("E38").Value = "  -1.67%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "16.80"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -3.04%  "

# Row 40
$ws.Range("B40").Value = "Stacks"
$ws.Range("C40").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.54"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -3.70%  "

# Row 41
$ws.Range("B41").Value = "Stellar"
$ws.Range("C41").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.115"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.06%  "

# Row 42
$ws.Range("E42").Value = "  -2.71%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "121.59"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.24%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "21.31"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -2.24%  "

# Row 46
$ws.Range("E46").Value = "  -1.85%  "

# Row 47
$ws.Range("E47").Value = "  +2.56%  "

# Row 48
$ws.Range("E48").Value = "  +0.90%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.011.05"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.83%  "

# Row 50
$ws.Range("E50").Value = "  +6.69%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.12"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +14.70%  "
